$wb = $excel.ActiveWorkbook

# Insert a new worksheet "total_concentrations" right after "input_concentrations"
# (i.e. right before "equilibrium_concentrations").
$afterSheet = $wb.Worksheets.Item("input_concentrations")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "total_concentrations"

# Header row
$newSheet.Range("A1").Value = "molecule1"
$newSheet.Range("B1").Value = "molecule2"
$newSheet.Range("C1").Value = "molecule3"

# Data rows
$newSheet.Range("A2").Value = 0.01
$newSheet.Range("B2").Value = 0.02
$newSheet.Range("C2").Value = 0.01

$newSheet.Range("A3").Value = 0.001
$newSheet.Range("B3").Value = 0.02
$newSheet.Range("C3").Value = 0.01

$newSheet.Range("A4").Value = 0.0001
$newSheet.Range("B4").Value = 0.01
$newSheet.Range("C4").Value = 0.005

# Restore the originally active sheet/tab (inserting a sheet makes the new
# sheet active by default).
$wb.Worksheets.Item("input_stoich_coefficients").Activate()
